# Add a team record (Wins/Losses/Ties) to the SEA_2007 sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row: AD1, AE1, AF1 -- match the style used for the other headers (A1:AC1)
$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"
$ws.Range("A1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122)  # xlPasteFormats

# Data rows 2-46: constant team record for every player row
$lastRow = 46
for ($r = 2; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, 30).Value = 88  # AD
    $ws.Cells.Item($r, 31).Value = 74  # AE
    $ws.Cells.Item($r, 32).Value = 0   # AF
}
